$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Locate the target paragraph: "Architecture centralisée/décentralisée" (lvl 2, 8th paragraph)
$target = $tr.Paragraphs(8, 1)

# Insert a new paragraph before it: "Prise en charge du réseau"
[void]$target.InsertBefore("Prise en charge du réseau`r")

# The new paragraph is now #8, the original (untouched) paragraph shifted to #9.
$newPara = $tr.Paragraphs(8, 1)
$oldPara = $tr.Paragraphs(9, 1)

# Split "Prise en charge du réseau" into two runs: "Prise en charge " + "du réseau"
$splitAt = 16
$tailNew = $tr.Characters($newPara.Start + $splitAt, $newPara.Length - $splitAt - 1)
$tailNew.Text = $tailNew.Text

# Split "Architecture centralisée/décentralisée" into two runs: "Architecture " + "centralisée/décentralisée"
$splitAt2 = 13
$tailOld = $tr.Characters($oldPara.Start + $splitAt2, $oldPara.Length - $splitAt2 - 1)
$tailOld.Text = $tailOld.Text
